# Quarterly indexing esoteric bug-fix operation
#
# Column A holds date-serial values marking the start of each forecast
# quarter. The fix shifts every one of these dates forward: from the
# 1st of its month to the 15th of the *following* month (i.e. the
# quarter-start index was off, and is corrected to land on the 15th of
# the next month instead of the 1st of the current one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $oldVal = $cell.Value2

    if ($oldVal -eq $null) { continue }

    $d = [DateTime]::FromOADate($oldVal)
    $shifted = $d.AddMonths(1)
    $newDate = Get-Date -Year $shifted.Year -Month $shifted.Month -Day 15 -Hour 0 -Minute 0 -Second 0

    $cell.Value = $newDate.ToOADate()
}
